# This script updates the "想去人数" (number of people who want to go)
# column (F) values on the "展览" and "全部类型" worksheets, matching the
# values published in the regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 377
    3  = 10865
    5  = 982
    6  = 188
    7  = 1345
    8  = 8321
    9  = 43
    10 = 469
    11 = 615
    12 = 222
    13 = 138
    14 = 3324
    17 = 39
    18 = 816
    22 = 128
    23 = 1825
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
